$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '24.825.87'
$ws.Range("E2").Value = '  +2.10%  '
$ws.Range("D3").Value = '1.664.98'
$ws.Range("E3").Value = '  +2.07%  '
$ws.Range("E4").Value = '  -0.49%  '
$ws.Range("D5").Value = '''330.24'
$ws.Range("E5").Value = '  +8.60%  '
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").Value = '''0.3646'
$ws.Range("E7").Value = '  +0.98%  '
$ws.Range("D8").Value = '''47.33'
$ws.Range("E8").Value = '  +0.98%  '
$ws.Range("D9").Value = '''0.3256'
$ws.Range("E9").Value = '  +1.13%  '
$ws.Range("D10").Value = '''1.137'
$ws.Range("E10").Value = '  +3.63%  '
$ws.Range("D11").Value = '''0.07067'
$ws.Range("E11").Value = '  +3.04%  '
$ws.Range("D12").Value = '''1.003'
$ws.Range("E12").Value = '  -0.01%  '
$ws.Range("D13").Value = '''6.068'
$ws.Range("E13").Value = '  +2.99%  '
$ws.Range("D14").Value = '''19.50'
$ws.Range("E14").Value = '  +2.46%  '
$ws.Range("D15").Value = '1.665.52'
$ws.Range("E15").Value = '  +1.83%  '
$ws.Range("D16").Value = '''6.587'
$ws.Range("E16").Value = '  +1.33%  '
$ws.Range("D17").Value = '''0.00001048'
$ws.Range("E17").Value = '  +0.93%  '
$ws.Range("D18").Value = '''0.06640'
$ws.Range("E18").Value = '  +2.10%  '
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("D20").Value = '''78.45'
$ws.Range("E20").Value = '  +3.14%  '
$ws.Range("D21").Value = '''5.917'
$ws.Range("E21").Value = '  +1.14%  '
$ws.Range("E22").Value = '  +0.85%  '
$ws.Range("D23").Value = '''12.52'
$ws.Range("E23").Value = '  +5.00%  '
$ws.Range("D24").Value = '24.793.23'
$ws.Range("E24").Value = '  +1.98%  '
$ws.Range("D25").Value = '''2.472'
$ws.Range("E25").Value = '  +3.42%  '
$ws.Range("D26").Value = '''2.424'
$ws.Range("E26").Value = '  +5.18%  '
$ws.Range("D27").Value = '''148.89'
$ws.Range("E27").Value = '  +3.68%  '
$ws.Range("D28").Value = '''18.64'
$ws.Range("E28").Value = '  +1.06%  '
$ws.Range("D29").Value = '1.849.20'
$ws.Range("E29").Value = '  +1.73%  '
$ws.Range("D30").Value = '''125.93'
$ws.Range("E30").Value = '  +1.78%  '
$ws.Range("D31").Value = '''1.167'
$ws.Range("E31").Value = '  +5.79%  '
$ws.Range("D32").Value = '''4.070'
$ws.Range("E32").Value = '  -0.11%  '
$ws.Range("D33").Value = '''5.700'
$ws.Range("E33").Value = '  +2.10%  '
$ws.Range("D34").Value = '''0.08494'
$ws.Range("E34").Value = '  +1.78%  '
$ws.Range("D35").Value = '''1.645'
$ws.Range("E35").Value = '  -1.15%  '
$ws.Range("D36").Value = '''12.16'
$ws.Range("E36").Value = '  -0.36%  '
$ws.Range("D37").Value = '''0.06193'
$ws.Range("E37").Value = '  +3.73%  '
$ws.Range("D38").Value = '''5.166'
$ws.Range("E38").Value = '  +1.77%  '
$ws.Range("D39").Value = '''0.02279'
$ws.Range("E39").Value = '  +3.57%  '
$ws.Range("E40").Value = '  +3.92%  '
$ws.Range("D41").Value = '''0.2083'
$ws.Range("E41").Value = '  +3.09%  '
$ws.Range("D42").Value = '''8.230'
$ws.Range("E42").Value = '  +1.73%  '
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").Value = '''0.5913'
$ws.Range("E44").Value = '  +1.98%  '
$ws.Range("D45").Value = '''3.843'
$ws.Range("E45").Value = '  +3.82%  '
$ws.Range("D46").Value = '''13.31'
$ws.Range("E46").Value = '  +7.17%  '
$ws.Range("D47").Value = '''0.5669'
$ws.Range("E47").Value = '  +3.34%  '
$ws.Range("D48").Value = '''125.62'
$ws.Range("E48").Value = '  +3.80%  '
$ws.Range("D49").Value = '''1.944'
$ws.Range("E49").Value = '  +2.02%  '
$ws.Range("D50").Value = '''0.06973'
$ws.Range("E50").Value = '  +1.65%  '
$ws.Range("D51").Value = '''1.191'
$ws.Range("E51").Value = '  +4.87%  '
